# December Deskcount update:
#  - Flip "Include in Occupancy Calculation" (column F) from Yes to No
#    for a handful of office rows.
#  - Bump the Melbourne headcount (column C, row 44) from 30 to 32.
#  - Leave the sheet scrolled/selected near the bottom, matching where
#    the author finished editing.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Deskcount")
$ws.Activate()

# Rows whose "Include in Occupancy Calculation" flag changes to "No".
$noRows = 16, 37, 38, 47, 48, 50
foreach ($r in $noRows) {
    $ws.Cells.Item($r, 6).Value = "No"
}

# Melbourne (row 44) headcount correction.
$ws.Cells.Item(44, 3).Value = 32

# Match the saved selection/scroll position.
$excel.ActiveWindow.ScrollRow = 28
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D57").Select()
